$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$textCells = @("D5","D6","D8","D12","D13","D14","D19","D20","D21","D22","D23","D24","D26","D27","D28","D31","D32","D33","D35","D36","D38","D39","D42","D43","D45","D46","D47","D48","D49","D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.369.64'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '2.504.47'
$ws.Range('E3').Value = '  -5.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '581.75'
$ws.Range('E5').Value = '  -2.29%  '
$ws.Range('D6').Value = '170.42'
$ws.Range('E6').Value = '  +1.24%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.524'
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('D9').Value = '2.503.55'
$ws.Range('E9').Value = '  -5.04%  '
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').Value = '5.11'
$ws.Range('E12').Value = '  -2.42%  '
$ws.Range('B13').Value = 'Cardano'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D13').Value = '0.348'
$ws.Range('E13').Value = '  -4.90%  '
$ws.Range('D14').Value = '26.63'
$ws.Range('E14').Value = '  -3.76%  '
$ws.Range('D15').Value = '2.957.09'
$ws.Range('E16').Value = '  -3.44%  '
$ws.Range('D17').Value = '66.258.43'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '2.512.35'
$ws.Range('E18').Value = '  -5.00%  '
$ws.Range('D19').Value = '11.23'
$ws.Range('E19').Value = '  -7.13%  '
$ws.Range('D20').Value = '7.69'
$ws.Range('E20').Value = '  -4.71%  '
$ws.Range('D21').Value = '346.88'
$ws.Range('E21').Value = '  -3.18%  '
$ws.Range('D22').Value = '4.18'
$ws.Range('E22').Value = '  -3.72%  '
$ws.Range('D23').Value = '4.62'
$ws.Range('E23').Value = '  -1.86%  '
$ws.Range('D24').Value = '1.96'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '69.64'
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').Value = '9.92'
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('D29').Value = '2.631.51'
$ws.Range('E29').Value = '  -4.75%  '
$ws.Range('D30').Value = '0.0₃0972'
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('D31').Value = '521.79'
$ws.Range('E31').Value = '  -5.09%  '
$ws.Range('D32').Value = '8.06'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').Value = '1.31'
$ws.Range('E33').Value = '  -3.68%  '
$ws.Range('E34').Value = '  -3.75%  '
$ws.Range('D35').Value = '0.130'
$ws.Range('E35').Value = '  -5.08%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  -3.54%  '
$ws.Range('D38').Value = '156.56'
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('D39').Value = '18.57'
$ws.Range('E39').Value = '  -2.56%  '
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').Value = '5.07'
$ws.Range('E43').Value = '  -3.48%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = '2.49'
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('D46').Value = '39.34'
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('D47').Value = '148.32'
$ws.Range('E47').Value = '  -3.27%  '
$ws.Range('D48').Value = '0.556'
$ws.Range('E48').Value = '  -4.66%  '
$ws.Range('D49').Value = '3.67'
$ws.Range('E49').Value = '  -4.10%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').Value = '1.71'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0270'
$ws.Range('E51').Value = '  -10.54%  '
